$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Intermediari_Partner")

$ws.Cells.Item(4, 3).Value = 759.0
$ws.Cells.Item(5, 3).Value = 492.0
$ws.Cells.Item(6, 3).Value = 442.0
$ws.Cells.Item(9, 3).Value = 386.0
$ws.Cells.Item(12, 3).Value = 325.0
$ws.Cells.Item(13, 3).Value = 235.0
$ws.Cells.Item(14, 1).Value = 'Progetti e Soluzioni SPA'
$ws.Cells.Item(14, 2).NumberFormat = "@"
$ws.Cells.Item(14, 2).Value = '06423240727'
$ws.Cells.Item(14, 3).Value = 181.0
$ws.Cells.Item(15, 1).Value = 'Alto Adige Riscossioni Spa'
$ws.Cells.Item(15, 2).NumberFormat = "@"
$ws.Cells.Item(15, 2).Value = '02805390214'
$ws.Cells.Item(15, 3).Value = 180.0
$ws.Cells.Item(16, 3).Value = 154.0
$ws.Cells.Item(17, 3).Value = 148.0
$ws.Cells.Item(18, 3).Value = 144.0
$ws.Cells.Item(19, 3).Value = 133.0
$ws.Cells.Item(20, 1).Value = 'Regione Piemonte'
$ws.Cells.Item(20, 2).NumberFormat = "@"
$ws.Cells.Item(20, 2).Value = '80087670016'
$ws.Cells.Item(20, 3).Value = 121.0
$ws.Cells.Item(21, 1).Value = 'Regione Puglia'
$ws.Cells.Item(21, 2).NumberFormat = "@"
$ws.Cells.Item(21, 2).Value = '80017210727'
$ws.Cells.Item(21, 3).Value = 120.0
$ws.Cells.Item(22, 1).Value = 'Regione Basilicata'
$ws.Cells.Item(22, 2).NumberFormat = "@"
$ws.Cells.Item(22, 2).Value = '80002950766'
$ws.Cells.Item(22, 3).Value = 106.0
$ws.Cells.Item(23, 1).Value = 'Regione Toscana'
$ws.Cells.Item(23, 2).NumberFormat = "@"
$ws.Cells.Item(23, 2).Value = '01386030488'
$ws.Cells.Item(23, 3).Value = 103.0
$ws.Cells.Item(25, 3).Value = 82.0
$ws.Cells.Item(27, 1).Value = 'Regione Autonoma della Sardegna'
$ws.Cells.Item(27, 2).NumberFormat = "@"
$ws.Cells.Item(27, 2).Value = '80002870923'
$ws.Cells.Item(27, 3).Value = 69.0
$ws.Cells.Item(28, 1).Value = 'Intesa Sanpaolo SPA'
$ws.Cells.Item(28, 2).NumberFormat = "@"
$ws.Cells.Item(28, 2).Value = '00799960158'
$ws.Cells.Item(28, 3).Value = 61.0
$ws.Cells.Item(29, 1).Value = 'PMPay s.r.l.'
$ws.Cells.Item(29, 2).NumberFormat = "@"
$ws.Cells.Item(29, 2).Value = '08747230962'
$ws.Cells.Item(29, 3).Value = 56.0
$ws.Cells.Item(30, 1).Value = 'NORDCOM'
$ws.Cells.Item(30, 2).NumberFormat = "@"
$ws.Cells.Item(30, 2).Value = '13384100155'
$ws.Cells.Item(30, 3).Value = 55.0
$ws.Cells.Item(31, 1).Value = 'Regione Umbria'
$ws.Cells.Item(31, 2).NumberFormat = "@"
$ws.Cells.Item(31, 2).Value = '80000130544'
$ws.Cells.Item(31, 3).Value = 52.0
$ws.Cells.Item(32, 1).Value = 'Bluenext S.r.l.'
$ws.Cells.Item(32, 2).NumberFormat = "@"
$ws.Cells.Item(32, 2).Value = '04228480408'
$ws.Cells.Item(32, 3).Value = 50.0
$ws.Cells.Item(33, 1).Value = 'Italriscossioni Società Italiana di Fiscalità Locale S.r.l.'
$ws.Cells.Item(33, 2).NumberFormat = "@"
$ws.Cells.Item(33, 2).Value = '06092371001'
$ws.Cells.Item(33, 3).Value = 46.0
$ws.Cells.Item(34, 1).Value = 'CINECA consorzio universitario'
$ws.Cells.Item(34, 2).NumberFormat = "@"
$ws.Cells.Item(34, 2).Value = '00317740371'
$ws.Cells.Item(34, 3).Value = 42.0
$ws.Cells.Item(35, 1).Value = 'Consorzio I.T. Srl'
$ws.Cells.Item(35, 2).NumberFormat = "@"
$ws.Cells.Item(35, 2).Value = '01321400192'
$ws.Cells.Item(35, 3).Value = 40.0
$ws.Cells.Item(36, 1).Value = 'ROMA CAPITALE'
$ws.Cells.Item(36, 2).NumberFormat = "@"
$ws.Cells.Item(36, 2).Value = '02438750586'
$ws.Cells.Item(36, 3).Value = 39.0
$ws.Cells.Item(37, 1).Value = 'UNIMATICA S.P.A'
$ws.Cells.Item(37, 2).NumberFormat = "@"
$ws.Cells.Item(37, 2).Value = '02098391200'
$ws.Cells.Item(40, 3).Value = 27.0
$ws.Cells.Item(41, 1).Value = 'Siscom SPA'
$ws.Cells.Item(41, 2).NumberFormat = "@"
$ws.Cells.Item(41, 2).Value = '01778000040'
$ws.Cells.Item(41, 3).Value = 26.0
$ws.Cells.Item(42, 1).Value = 'Regione Liguria'
$ws.Cells.Item(42, 2).NumberFormat = "@"
$ws.Cells.Item(42, 2).Value = '00849050109'
$ws.Cells.Item(43, 1).Value = 'Novares Spa'
$ws.Cells.Item(43, 2).NumberFormat = "@"
$ws.Cells.Item(43, 2).Value = '12105121003'
$ws.Cells.Item(43, 3).Value = 19.0
$ws.Cells.Item(44, 1).Value = 'Citta'' Metropolitana di Roma Capitale'
$ws.Cells.Item(44, 2).NumberFormat = "@"
$ws.Cells.Item(44, 2).Value = '80034390585'
$ws.Cells.Item(45, 1).Value = 'Nexi SpA'
$ws.Cells.Item(45, 2).NumberFormat = "@"
$ws.Cells.Item(45, 2).Value = '13212880150'
$ws.Cells.Item(46, 1).Value = 'Regione Lazio'
$ws.Cells.Item(46, 2).NumberFormat = "@"
$ws.Cells.Item(46, 2).Value = '80143490581'
$ws.Cells.Item(47, 1).Value = 'ANDREANI TRIBUTI srl'
$ws.Cells.Item(47, 2).NumberFormat = "@"
$ws.Cells.Item(47, 2).Value = '01412920439'
$ws.Cells.Item(47, 3).Value = 18.0
$ws.Cells.Item(48, 1).Value = 'Comune di Palermo'
$ws.Cells.Item(48, 2).NumberFormat = "@"
$ws.Cells.Item(48, 2).Value = '80016350821'
$ws.Cells.Item(49, 1).Value = 'Numera Sistemi e Informatica SpA'
$ws.Cells.Item(49, 2).NumberFormat = "@"
$ws.Cells.Item(49, 2).Value = '01265230902'
$ws.Cells.Item(49, 3).Value = 17.0
$ws.Cells.Item(50, 1).Value = 'Servizi Locali SpA'
$ws.Cells.Item(50, 2).NumberFormat = "@"
$ws.Cells.Item(50, 2).Value = '03170580751'
$ws.Cells.Item(50, 3).Value = 14.0
$ws.Cells.Item(51, 1).Value = 'Regione Autonoma Valle D''Aosta'
$ws.Cells.Item(51, 2).NumberFormat = "@"
$ws.Cells.Item(51, 2).Value = '80002270074'
$ws.Cells.Item(51, 3).Value = 13.0
$ws.Cells.Item(58, 1).Value = 'ARCA Servizi s.r.l'
$ws.Cells.Item(58, 2).NumberFormat = "@"
$ws.Cells.Item(58, 2).Value = '09106071005'
$ws.Cells.Item(58, 3).Value = 4.0
$ws.Cells.Item(59, 1).Value = 'Linea Comune Spa'
$ws.Cells.Item(59, 2).NumberFormat = "@"
$ws.Cells.Item(59, 2).Value = '05591950489'
$ws.Cells.Item(61, 1).Value = 'ISWEB S.p.A.'
$ws.Cells.Item(61, 2).NumberFormat = "@"
$ws.Cells.Item(61, 2).Value = '01722270665'
$ws.Cells.Item(61, 3).Value = 3.0
$ws.Cells.Item(62, 1).Value = 'CityPoste Payment Digital S.r.l.'
$ws.Cells.Item(62, 2).NumberFormat = "@"
$ws.Cells.Item(62, 2).Value = '02003750672'
$ws.Cells.Item(63, 1).Value = 'ICCREA Banca SpA'
$ws.Cells.Item(63, 2).NumberFormat = "@"
$ws.Cells.Item(63, 2).Value = '04774801007'
$ws.Cells.Item(64, 1).Value = 'Agenzia Italiana del Farmaco - AIFA'
$ws.Cells.Item(64, 2).NumberFormat = "@"
$ws.Cells.Item(64, 2).Value = '97345810580'
$ws.Cells.Item(65, 1).Value = 'MegASP S.r.l.'
$ws.Cells.Item(65, 2).NumberFormat = "@"
$ws.Cells.Item(65, 2).Value = '09898030151'
$ws.Cells.Item(66, 1).Value = 'Ministero dello Sviluppo Economico'
$ws.Cells.Item(66, 2).NumberFormat = "@"
$ws.Cells.Item(66, 2).Value = '80230390587'
$ws.Cells.Item(67, 1).Value = 'Softline srl'
$ws.Cells.Item(67, 2).NumberFormat = "@"
$ws.Cells.Item(67, 2).Value = '12299030150'
$ws.Cells.Item(68, 1).Value = 'I.C.A. - Imposte Comunali Affini – s.r.l.'
$ws.Cells.Item(68, 2).NumberFormat = "@"
$ws.Cells.Item(68, 2).Value = '02478610583'
$ws.Cells.Item(69, 1).Value = 'Banco BPM Società per Azioni'
$ws.Cells.Item(69, 2).NumberFormat = "@"
$ws.Cells.Item(69, 2).Value = '09722490969'
$ws.Cells.Item(70, 1).Value = 'ARGO SOFTWARE SRL'
$ws.Cells.Item(70, 2).NumberFormat = "@"
$ws.Cells.Item(70, 2).Value = '00838520880'
$ws.Cells.Item(71, 1).Value = 'Engineering Ingegneria Informatica SpA'
$ws.Cells.Item(71, 2).NumberFormat = "@"
$ws.Cells.Item(71, 2).Value = '00967720285'
$ws.Cells.Item(72, 1).Value = 'BANCA MONTE DEI PASCHI DI SIENA'
$ws.Cells.Item(72, 2).NumberFormat = "@"
$ws.Cells.Item(72, 2).Value = '00884060526'
$ws.Cells.Item(73, 1).Value = 'San Marco SPA'
$ws.Cells.Item(73, 2).NumberFormat = "@"
$ws.Cells.Item(73, 2).Value = '04142440728'
